# "cập nhật mới nhất" - refresh the leave-tracking export with the latest
# HR data: reorder a couple of header columns and replace the data rows
# with the newest pull (6 employee rows instead of 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: Thai San / Viec Rieng swap ahead of Hieu Hi / Phep Nam ---
$ws.Range("K1").Value = 'Thai Sản'
$ws.Range("L1").Value = 'Việc Riêng'
$ws.Range("M1").Value = 'Hiếu Hỉ'
$ws.Range("N1").Value = 'Phép Năm'

# --- Row 2: new employee (Trương Văn Dũng) ---
# (I2,K2:P2 were already blank on this physical row, so they're left as-is)
$ws.Range("A2").Value = 621605001
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 2023
$ws.Range("D2").Value = 'DŨNG'
$ws.Range("E2").Value = 'TRƯƠNG VĂN'
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = '08/01/2016'
$ws.Range("G2").Value = 'TN Quản Trị Mạng BP CNTT'
$ws.Range("H2").Value = 'CNTT'
$ws.Range("J2").Value = 2

# --- Row 3: Lê Viết Hải (was row 2), refreshed totals ---
# (I3,M3,N3 held the old Sơn row's AnnualLeave/Thai San/Viec Rieng numbers
# on this physical row, so those need clearing; O3/P3 were already blank)
$ws.Range("A3").Value = 622204064
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 2023
$ws.Range("D3").Value = 'HẢI'
$ws.Range("E3").Value = 'LÊ VIẾT'
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = '08/01/2018'
$ws.Range("G3").Value = 'NV PT Phần Mềm BP CNTT'
$ws.Range("H3").Value = 'CNTT'
$ws.Range("I3").ClearContents()
$ws.Range("J3").Value = 9
$ws.Range("K3").Value = 10
$ws.Range("L3").Value = 7
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

# --- Row 4: Ngô Quốc Sơn (was row 3), refreshed totals ---
# (O4/P4 were already blank on this physical row, so left as-is)
$ws.Range("A4").Value = 622210012
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 2023
$ws.Range("D4").Value = 'SƠN'
$ws.Range("E4").Value = 'NGÔ QUỐC'
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '08/01/2020'
$ws.Range("G4").Value = 'NV PT Phần Mềm BP CNTT'
$ws.Range("H4").Value = 'CNTT'
$ws.Range("I4").Value = 40
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 18
$ws.Range("M4").Value = 3
$ws.Range("N4").Value = 2

# --- Row 5: new employee (Đoàn Trọng Nghĩa) ---
$ws.Range("A5").Value = 621805013
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 2023
$ws.Range("D5").Value = 'NGHĨA'
$ws.Range("E5").Value = 'ĐOÀN TRỌNG'
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = '05/18/2018'
$ws.Range("G5").Value = 'NV Hành chính BP HCNS'
$ws.Range("H5").Value = 'HCNS'
$ws.Range("L5").Value = 1

# --- Row 6: new employee (Phan Như Ý) ---
$ws.Range("A6").Value = 621412001
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 2023
$ws.Range("D6").Value = 'Ý'
$ws.Range("E6").Value = 'PHAN NHƯ'
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = '12/15/2014'
$ws.Range("G6").Value = 'NV CS Tiền lương BP HCNS'
$ws.Range("H6").Value = 'HCNS'
$ws.Range("L6").Value = 1
